$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for the Price/Volume columns so that numeric-looking
# strings (e.g. trailing zeros, thousand-dot separators) are not coerced to numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '41.687.50'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').Value = '2.479.21'
$ws.Range('E3').Value = '  -0.18%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '319.39'
$ws.Range('E5').Value = '  +1.57%  '
$ws.Range('D6').Value = '93.12'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').Value = '  +1.98%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.519'
$ws.Range('E9').Value = '  +2.30%  '
$ws.Range('D10').Value = '0.0891'
$ws.Range('E10').Value = '  +13.46%  '
$ws.Range('D11').Value = '33.08'
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('D13').Value = '2.856.86'
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('D14').Value = '6.97'
$ws.Range('E14').Value = '  +1.75%  '
$ws.Range('D15').Value = '15.71'
$ws.Range('E15').Value = '  -2.84%  '
$ws.Range('D16').Value = '2.469.32'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('D17').Value = '0.793'
$ws.Range('E17').Value = '  +3.90%  '
$ws.Range('D18').Value = '41.632.72'
$ws.Range('D19').Value = '0.0₃0964'
$ws.Range('E19').Value = '  +2.96%  '
$ws.Range('D20').Value = '6.50'
$ws.Range('E20').Value = '  +1.61%  '
$ws.Range('D21').Value = '71.55'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').Value = '11.52'
$ws.Range('E22').Value = '  +1.24%  '
$ws.Range('D23').Value = '241.59'
$ws.Range('E23').Value = '  +1.94%  '
$ws.Range('D24').Value = '2.76'
$ws.Range('E24').Value = '  +1.83%  '
$ws.Range('E25').Value = '  +1.55%  '
$ws.Range('D27').Value = '24.95'
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').Value = '2.30'
$ws.Range('E28').Value = '  +4.64%  '
$ws.Range('D29').Value = '9.90'
$ws.Range('E29').Value = '  +2.30%  '
$ws.Range('D30').Value = '36.56'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('D31').Value = '156.68'
$ws.Range('E31').Value = '  -0.63%  '
$ws.Range('D32').Value = '5.55'
$ws.Range('E32').Value = '  +0.85%  '
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('D34').Value = '0.0772'
$ws.Range('E34').Value = '  +2.63%  '
$ws.Range('D35').Value = '2.58'
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('D36').Value = '17.61'
$ws.Range('E36').Value = '  +0.90%  '
$ws.Range('D37').Value = '2.93'
$ws.Range('E37').Value = '  -0.27%  '
$ws.Range('D38').Value = '1.84'
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('E39').Value = '  +1.22%  '
$ws.Range('D40').Value = '0.104'
$ws.Range('E40').Value = '  -1.18%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '4.01'
$ws.Range('E41').Value = '  -2.35%  '
$ws.Range('B42').Value = 'ApeXProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D42').Value = '2.50'
$ws.Range('E42').Value = '  +1.82%  '
$ws.Range('D43').Value = '1.984.39'
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').Value = '19.42'
$ws.Range('E44').Value = '  -2.57%  '
$ws.Range('D45').Value = '0.0286'
$ws.Range('E45').Value = '  +0.73%  '
$ws.Range('D46').Value = '3.03'
$ws.Range('E46').Value = '  +2.71%  '
$ws.Range('D47').Value = '9.17'
$ws.Range('E47').Value = '  +0.90%  '
$ws.Range('D48').Value = '2.710.86'
$ws.Range('E48').Value = '  -0.46%  '
$ws.Range('D49').Value = '97.71'
$ws.Range('E49').Value = '  -0.28%  '
$ws.Range('D50').Value = '67.88'
$ws.Range('E50').Value = '  -0.26%  '
$ws.Range('D51').Value = '74.05'
$ws.Range('E51').Value = '  +2.20%  '

# Restore the original (default) cell style now that values are safely stored as text.
$dataRange.Style = "Normal"
